# This script applies the row-level Price/Volume(1h) refresh (and the
# Mantle/Hedera and VeChain/WhiteBITCoin row-order swaps) described by the
# target diff for cryptos.xlsx.
#
# Numeric-looking text values (e.g. "558.14", "1.00") are written with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cells) instead of converting them to numbers, and the cell
# style is reset to "Normal" afterwards so no stray number-format/quote-
# prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '63.314.19'
$ws.Range("E2").Value = '  -1.30%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '2.715.09'
$ws.Range("E3").Value = '  -1.71%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.06%  '
# Row 5: BNB
$ws.Range("D5").Value = '''558.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.42%  '
# Row 6: Solana
$ws.Range("D6").Value = '''157.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.61%  '
# Row 7: USDC
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
# Row 8: XRP
$ws.Range("D8").Value = '''0.593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.07%  '
# Row 9: Dogecoin
$ws.Range("E9").Value = '  -3.44%  '
# Row 10: TRON
$ws.Range("E10").Value = '  -0.06%  '
# Row 11: Toncoin
$ws.Range("D11").Value = '''5.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.99%  '
# Row 12: Cardano
$ws.Range("D12").Value = '''0.373'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.30%  '
# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.194.35'
$ws.Range("E13").Value = '  -1.74%  '
# Row 14: Avalanche
$ws.Range("D14").Value = '''26.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.58%  '
# Row 15: WrappedBTC
$ws.Range("D15").Value = '63.194.75'
$ws.Range("E15").Value = '  -0.89%  '
# Row 16: ShibaInu
$ws.Range("E16").Value = '  -4.15%  '
# Row 17: WrappedEther
$ws.Range("D17").Value = '2.714.61'
$ws.Range("E17").Value = '  -1.79%  '
# Row 18: Chainlink
$ws.Range("E18").Value = '  -0.93%  '
# Row 19: Polkadot
$ws.Range("D19").Value = '''4.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.87%  '
# Row 20: BitcoinCash
$ws.Range("D20").Value = '''350.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.58%  '
# Row 21: Uniswap
$ws.Range("D21").Value = '''6.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.40%  '
# Row 22: Dai
$ws.Range("D22").Value = '''1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.03%  '
# Row 23: Polygon
$ws.Range("D23").Value = '''0.513'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.68%  '
# Row 24: Litecoin
$ws.Range("D24").Value = '''64.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.17%  '
# Row 25: Kaspa
$ws.Range("E25").Value = '  -1.65%  '
# Row 26: Binance-PegBSC-USD
$ws.Range("E26").Value = '  +0.17%  '
# Row 27: InternetComputer(DFINITY)
$ws.Range("E27").Value = '  -5.42%  '
# Row 28: PEPE
$ws.Range("D28").Value = [string]::Concat("0.0", [char]0x2083, "0881")
$ws.Range("E28").Value = '  -4.36%  '
# Row 29: Fetch.AI
$ws.Range("E29").Value = '  +9.66%  '
# Row 31: Aptos
$ws.Range("D31").Value = '''7.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.08%  '
# Row 32: Monero
$ws.Range("D32").Value = '''166.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.70%  '
# Row 33: ImmutableX
$ws.Range("E33").Value = '  -1.41%  '
# Row 34: USDe
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '
# Row 35: EthereumClassic
$ws.Range("D35").Value = '''19.80'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.40%  '
# Row 36: NEARProtocol
$ws.Range("D36").Value = '''4.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.78%  '
# Row 37: Stacks
$ws.Range("E37").Value = '  -2.87%  '
# Row 38: Bittensor
$ws.Range("D38").Value = '''346.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.51%  '
# Row 39: SuiNetwork
$ws.Range("D39").Value = '''0.957'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.28%  '
# Row 40: RenderToken
$ws.Range("D40").Value = '''6.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.65%  '
# Row 41: Filecoin
$ws.Range("D41").Value = '''4.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.61%  '
# Row 42: OKB
$ws.Range("D42").Value = '''38.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.44%  '
# Row 43: InjectiveProtocol
$ws.Range("D43").Value = '''21.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.66%  '
# Row 44: EnergySwap
$ws.Range("E44").Value = '  -4.33%  '
# Row 45: Mantle
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.627'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '
# Row 46: Hedera
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0571'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.00%  '
# Row 47: FirstDigitalUSD
$ws.Range("E47").Value = '  -0.07%  '
# Row 48: Aave
$ws.Range("D48").Value = '''131.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.47%  '
# Row 49: Stellar
$ws.Range("D49").Value = '''0.0984'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.71%  '
# Row 50: WhiteBITCoin
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '''11.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.11%  '
# Row 51: VeChain
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '''0.0245'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.98%  '
